# Add a new "2021" data column (column J) to the table, mirroring column I's
# layout, and select L27 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that are pure border/separator rows: copy I's cell format into J so
#     the border carries all the way across the newly-added column. ---
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I27").Copy()
$ws.Range("J27").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# --- Header row: 2021 year label, formatted like the other year headers. ---
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("J4").Value = 2021

# --- Data rows: new 2021 figures. These take on the ambient row/column
#     default formatting (no explicit override), same as the source data. ---
$ws.Range("J5").Value = 5356.3
$ws.Range("J6").Value = 9.5
$ws.Range("J8").Value = 7.9
$ws.Range("J9").Value = 10.5
$ws.Range("J11").Value = 9.6
$ws.Range("J12").Value = 9.4
$ws.Range("J14").Value = 14.8
$ws.Range("J15").Value = 9.1
$ws.Range("J16").Value = 9.5
$ws.Range("J17").Value = 5.9

# --- Final selection, matching the author's last recorded cursor cell. ---
$ws.Range("L27").Select()
